$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 updates (becomes old row 13's data; AC11 comment is removed)
$ws.Range("A11").Value = 112243230
$ws.Range("B11").Value = 89553
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 1202
$ws.Range("F11").Value = "Ullticka"
$ws.Range("G11").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H11").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Z11").Value = "10:16"
$ws.Range("AB11").Value = "10:16"
$ws.Range("AC11").ClearContents()
$ws.Range("AJ11").Value = "gran"
$ws.Range("AK11").Value = "Picea abies"
$ws.Range("AO11").Value = "Picea abies"

# Row 12 updates (becomes old row 11's data; gains AC12 comment)
$ws.Range("A12").Value = 112243622
$ws.Range("B12").Value = 89893
$ws.Range("D12").Value = "VU"
$ws.Range("E12").Value = 2062
$ws.Range("F12").Value = "Ulltickeporing"
$ws.Range("G12").Value = "Skeletocutis brevispora"
$ws.Range("H12").Value = "Niemelä"
$ws.Range("Q12").Value = 523006
$ws.Range("R12").Value = 6739484
$ws.Range("Z12").Value = "10:21"
$ws.Range("AB12").Value = "10:21"
$ws.Range("AC12").Value = "På ytmurken granlåga med delvis avfallande bark och insektsgnag i veden från tiden då granen stod upp."
$ws.Range("AJ12").Value = "ullticka"
$ws.Range("AK12").Value = "Phellinidium ferrugineofuscum"
$ws.Range("AO12").Value = "Phellinidium ferrugineofuscum"

# Row 13 updates (becomes old row 12's data)
$ws.Range("A13").Value = 112244426
$ws.Range("B13").Value = 77685
$ws.Range("E13").Value = 185
$ws.Range("F13").Value = "Violettgrå tagellav"
$ws.Range("G13").Value = "Bryoria nadvornikiana"
$ws.Range("H13").Value = "(Gyeln.) Brodo & D.Hawksw."
$ws.Range("Q13").Value = 523094
$ws.Range("R13").Value = 6739613
$ws.Range("Z13").Value = "11:11"
$ws.Range("AB13").Value = "11:11"
